$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Neo4j query text that documents how the companion *_Neo4jData.xlsx
# file for this test case was produced. It goes into A2, next to the
# existing WebData/Neo4jData file-name cells in B2/C2.
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma of the colon''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Range("A2").Value = $query

# Row grows tall enough to show the (wrapped) query text.
$ws.Rows.Item(2).RowHeight = 87

# Leave the sheet scrolled/selected the way it was left after entering
# the query: A5:A10 highlighted.
$null = $ws.Range("A5:A10").Select()
